# EOH constraint working in monolithic model
# Raise the plateau value of each scenario row from 59 to 60, starting at the
# first column where the row's growth curve has flattened out (varies by row)
# through the last data column (BJ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 62  # column BJ

# Map of row number -> first column (1-based) that should become 60
$startCols = @{
    2  = 48  # AV
    3  = 49  # AW
    4  = 49  # AW
    5  = 50  # AX
    6  = 50  # AX
    7  = 50  # AX
    8  = 49  # AW
    9  = 50  # AX
    10 = 50  # AX
    11 = 51  # AY
    12 = 51  # AY
    13 = 52  # AZ
}

foreach ($row in $startCols.Keys) {
    $startCol = $startCols[$row]
    for ($col = $startCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = 60
    }
}
